$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7970000505447388
$ws.Range("C2").Value = 0.7885000705718994
$ws.Range("D2").Value = 0.796750009059906
$ws.Range("E2").Value = 0.8075000047683716
$ws.Range("F2").Value = 0.7775000333786011

# Row 3
$ws.Range("B3").Value = 0.8425000309944153
$ws.Range("C3").Value = 0.8415000438690186
$ws.Range("D3").Value = 0.8425000309944153
$ws.Range("E3").Value = 0.8472500443458557
$ws.Range("F3").Value = 0.8375000357627869

# Row 4
$ws.Range("B4").Value = 0.8010774850845337
$ws.Range("C4").Value = 0.7919024229049683
$ws.Range("D4").Value = 0.8207361698150635
$ws.Range("E4").Value = 0.8211656808853149
$ws.Range("F4").Value = 0.7915928959846497

# Row 5
$ws.Range("B5").Value = 0.8468078374862671
$ws.Range("C5").Value = 0.843259871006012
$ws.Range("D5").Value = 0.8492092490196228
$ws.Range("E5").Value = 0.8250664472579956
$ws.Range("F5").Value = 0.8230729103088379

# Row 6 (B6 newly added)
$ws.Range("B6").Value = 0.5290902853012085
$ws.Range("C6").Value = 0.7465657591819763
$ws.Range("D6").Value = 0.7491359710693359
$ws.Range("E6").Value = 0.8074912428855896
$ws.Range("F6").Value = 0.7802461385726929

# Row 7 (B7 newly added)
$ws.Range("B7").Value = 0.6302996873855591
$ws.Range("C7").Value = 0.8134551048278809
$ws.Range("D7").Value = 0.82331383228302
$ws.Range("E7").Value = 0.8510632514953613
$ws.Range("F7").Value = 0.8392308950424194

# Row 8
$ws.Range("B8").Value = 0.8135000467300415
$ws.Range("C8").Value = 0.6082500219345093
$ws.Range("D8").Value = 0.5239999890327454
$ws.Range("E8").Value = 0.8652500510215759
$ws.Range("F8").Value = 0.9072500467300415

# Row 9
$ws.Range("B9").Value = 0.01575000025331974
$ws.Range("C9").Value = 0.0755000039935112
$ws.Range("D9").Value = 0.3634999990463257
$ws.Range("E9").Value = 0.8657500743865967
$ws.Range("F9").Value = 0.9065000414848328

# Row 10
$ws.Range("B10").Value = 0.7903858423233032
$ws.Range("C10").Value = 0.7779074311256409
$ws.Range("D10").Value = 0.7889476418495178
$ws.Range("E10").Value = 0.8201858997344971
$ws.Range("F10").Value = 0.826413094997406

# Row 11 (E11 and F11 removed/cleared)
$ws.Range("B11").Value = 0.8331082463264465
$ws.Range("C11").Value = 0.8342833518981934
$ws.Range("D11").Value = 0.8491634726524353
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()

# Row 12
$ws.Range("B12").Value = 0.790385901927948
$ws.Range("C12").Value = 0.7779074311256409
$ws.Range("D12").Value = 0.7889476418495178
$ws.Range("E12").Value = 0.8201858997344971
$ws.Range("F12").Value = 0.826413094997406
